# Strike through the "b. Las operaciones son static, no se subrayaron." note
# (item 3.b in the "Corrección" list) to mark it as resolved, matching its
# sibling items 3.a and 3.c which are already struck through.

$d = $word.ActiveDocument

$target = "Las operaciones son static"
$count = 0

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*$target*") {
        $p.Range.Font.StrikeThrough = 1
        $count = $count + 1
    }
}

Write-Host "Paragraphs struck through: $count"
